$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update BOM row 2 (U$1) to reflect the new regulator part:
#   old: 296-20796-1-ND / IC REG LDO 5V 1.5A DDPAK @ 0.93
#   new: 945-1610-5-ND  /  CONV DC/DC 1.5A 6.5-18VIN 5V @ 10.73
$ws.Range("B2").Value = "945-1610-5-ND"
$ws.Range("C2").Value = " CONV DC/DC 1.5A 6.5-18VIN 5V"
$ws.Range("E2").Value = 10.73
$ws.Range("F2").Value = 10.73

# Move the active selection to F5 (matches the saved view state)
$ws.Range("F5").Select()
